# Fruta / hortaliza, semanal
# Insert a new weekly data record at row 463, shifting all subsequent
# rows (old 463..537) down by one (to 464..538).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 463, pushing existing rows down.
$ws.Rows.Item(463).Insert()

# Populate the newly inserted row 463 with the new record's data.
$ws.Range("A463").Value = 9
$ws.Range("B463").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C463").Value = "Metropolitana"
$ws.Range("D463").Value = 44984
$ws.Range("E463").Value = 13
$ws.Range("F463").Value = 100112032
$ws.Range("G463").Value = "Zapallo italiano"
$ws.Range("H463").Value = "Sin especificar"
$ws.Range("I463").Value = "Primera"
$ws.Range("J463").Value = 250
$ws.Range("K463").Value = 6000
$ws.Range("L463").Value = 7000
$ws.Range("M463").Value = 6500
$ws.Range("N463").Value = '$/caja 50 unidades'
$ws.Range("O463").Value = "Región de O'Higgins"
$ws.Range("P463").Value = 130
$ws.Range("Q463").Value = 50
$ws.Range("R463").Value = "Hortaliza"
